$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = "@"
$ws.Range('D2').Value = '63.508.46'
$ws.Range('E2').NumberFormat = "@"
$ws.Range('E2').Value = '  +7.06%  '
$ws.Range('D3').NumberFormat = "@"
$ws.Range('D3').Value = '3.118.53'
$ws.Range('E3').NumberFormat = "@"
$ws.Range('E3').Value = '  +4.64%  '
$ws.Range('E4').NumberFormat = "@"
$ws.Range('E4').Value = '  -0.04%  '
$ws.Range('D5').NumberFormat = "@"
$ws.Range('D5').Value = '586.42'
$ws.Range('E5').NumberFormat = "@"
$ws.Range('E5').Value = '  +3.69%  '
$ws.Range('D6').NumberFormat = "@"
$ws.Range('D6').Value = '145.45'
$ws.Range('E6').NumberFormat = "@"
$ws.Range('E6').Value = '  +5.69%  '
$ws.Range('E7').NumberFormat = "@"
$ws.Range('E7').Value = '  -0.01%  '
$ws.Range('D8').NumberFormat = "@"
$ws.Range('D8').Value = '3.111.33'
$ws.Range('E8').NumberFormat = "@"
$ws.Range('E8').Value = '  +4.72%  '
$ws.Range('D9').NumberFormat = "@"
$ws.Range('D9').Value = '0.532'
$ws.Range('E9').NumberFormat = "@"
$ws.Range('E9').Value = '  +2.48%  '
$ws.Range('E10').NumberFormat = "@"
$ws.Range('E10').Value = '  +14.91%  '
$ws.Range('E11').NumberFormat = "@"
$ws.Range('E11').Value = '  +8.31%  '
$ws.Range('D12').NumberFormat = "@"
$ws.Range('D12').Value = '0.469'
$ws.Range('E12').NumberFormat = "@"
$ws.Range('E12').Value = '  +4.32%  '
$ws.Range('D13').NumberFormat = "@"
$ws.Range('D13').Value = '0.0000248'
$ws.Range('E13').NumberFormat = "@"
$ws.Range('E13').Value = '  +8.60%  '
$ws.Range('D14').NumberFormat = "@"
$ws.Range('D14').Value = '35.54'
$ws.Range('E14').NumberFormat = "@"
$ws.Range('E14').Value = '  +5.82%  '
$ws.Range('E15').NumberFormat = "@"
$ws.Range('E15').Value = '  +0.86%  '
$ws.Range('D16').NumberFormat = "@"
$ws.Range('D16').Value = '3.633.58'
$ws.Range('E16').NumberFormat = "@"
$ws.Range('E16').Value = '  +4.63%  '
$ws.Range('E17').NumberFormat = "@"
$ws.Range('E17').Value = '  +1.50%  '
$ws.Range('D18').NumberFormat = "@"
$ws.Range('D18').Value = '63.372.27'
$ws.Range('E18').NumberFormat = "@"
$ws.Range('E18').Value = '  +6.79%  '
$ws.Range('D19').NumberFormat = "@"
$ws.Range('D19').Value = '3.111.81'
$ws.Range('E19').NumberFormat = "@"
$ws.Range('E19').Value = '  +4.37%  '
$ws.Range('D20').NumberFormat = "@"
$ws.Range('D20').Value = '468.50'
$ws.Range('E20').NumberFormat = "@"
$ws.Range('E20').Value = '  +7.44%  '
$ws.Range('D21').NumberFormat = "@"
$ws.Range('D21').Value = '14.15'
$ws.Range('E21').NumberFormat = "@"
$ws.Range('E21').Value = '  +4.41%  '
$ws.Range('D22').NumberFormat = "@"
$ws.Range('D22').Value = '0.728'
$ws.Range('E22').NumberFormat = "@"
$ws.Range('E22').Value = '  +1.20%  '
$ws.Range('E23').NumberFormat = "@"
$ws.Range('E23').Value = '  +7.77%  '
$ws.Range('E24').NumberFormat = "@"
$ws.Range('E24').Value = '  +1.00%  '
$ws.Range('E25').NumberFormat = "@"
$ws.Range('E25').Value = '  +2.62%  '
$ws.Range('E26').NumberFormat = "@"
$ws.Range('E26').Value = '  +0.04%  '
$ws.Range('D27').NumberFormat = "@"
$ws.Range('D27').Value = '8.66'
$ws.Range('E27').NumberFormat = "@"
$ws.Range('E27').Value = '  +12.46%  '
$ws.Range('E28').NumberFormat = "@"
$ws.Range('E28').Value = '  +0.53%  '
$ws.Range('E29').NumberFormat = "@"
$ws.Range('E29').Value = '  +5.04%  '
$ws.Range('E30').NumberFormat = "@"
$ws.Range('E30').Value = '  +0.00%  '
$ws.Range('D31').NumberFormat = "@"
$ws.Range('D31').Value = '6.84'
$ws.Range('E31').NumberFormat = "@"
$ws.Range('E31').Value = '  +9.90%  '
$ws.Range('D32').NumberFormat = "@"
$ws.Range('D32').Value = '27.04'
$ws.Range('E32').NumberFormat = "@"
$ws.Range('E32').Value = '  +5.20%  '
$ws.Range('E33').NumberFormat = "@"
$ws.Range('E33').Value = '  +5.55%  '
$ws.Range('D34').NumberFormat = "@"
$ws.Range('D34').Value = '0.0₃0875'
$ws.Range('E34').NumberFormat = "@"
$ws.Range('E34').Value = '  +14.20%  '
$ws.Range('E35').NumberFormat = "@"
$ws.Range('E35').Value = '  +17.09%  '
$ws.Range('E36').NumberFormat = "@"
$ws.Range('E36').Value = '  +7.08%  '
$ws.Range('B37').NumberFormat = "@"
$ws.Range('B37').Value = 'Filecoin'
$ws.Range('C37').NumberFormat = "@"
$ws.Range('C37').Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range('D37').NumberFormat = "@"
$ws.Range('D37').Value = '6.10'
$ws.Range('E37').NumberFormat = "@"
$ws.Range('E37').Value = '  +3.80%  '
$ws.Range('B38').NumberFormat = "@"
$ws.Range('B38').Value = 'dogwifhat'
$ws.Range('C38').NumberFormat = "@"
$ws.Range('C38').Value = 'https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif'
$ws.Range('D38').NumberFormat = "@"
$ws.Range('D38').Value = '3.34'
$ws.Range('E38').NumberFormat = "@"
$ws.Range('E38').Value = '  +21.39%  '
$ws.Range('D39').NumberFormat = "@"
$ws.Range('D39').Value = '50.72'
$ws.Range('E39').NumberFormat = "@"
$ws.Range('E39').Value = '  +4.41%  '
$ws.Range('D40').NumberFormat = "@"
$ws.Range('D40').Value = '443.98'
$ws.Range('E40').NumberFormat = "@"
$ws.Range('E40').Value = '  +11.21%  '
$ws.Range('E41').NumberFormat = "@"
$ws.Range('E41').Value = '  +0.29%  '
$ws.Range('D42').NumberFormat = "@"
$ws.Range('D42').Value = '2.911.97'
$ws.Range('E42').NumberFormat = "@"
$ws.Range('E42').Value = '  +6.74%  '
$ws.Range('E43').NumberFormat = "@"
$ws.Range('E43').Value = '  +5.31%  '
$ws.Range('D44').NumberFormat = "@"
$ws.Range('D44').Value = '0.280'
$ws.Range('E44').NumberFormat = "@"
$ws.Range('E44').Value = '  +12.13%  '
$ws.Range('E45').NumberFormat = "@"
$ws.Range('E45').Value = '  +5.29%  '
$ws.Range('E46').NumberFormat = "@"
$ws.Range('E46').Value = '  +9.12%  '
$ws.Range('D47').NumberFormat = "@"
$ws.Range('D47').Value = '35.96'
$ws.Range('E47').NumberFormat = "@"
$ws.Range('E47').Value = '  +5.06%  '
$ws.Range('D49').NumberFormat = "@"
$ws.Range('D49').Value = '123.72'
$ws.Range('E49').NumberFormat = "@"
$ws.Range('E49').Value = '  +1.31%  '
$ws.Range('E50').NumberFormat = "@"
$ws.Range('E50').Value = '  +1.49%  '
$ws.Range('D51').NumberFormat = "@"
$ws.Range('D51').Value = '24.67'
$ws.Range('E51').NumberFormat = "@"
$ws.Range('E51').Value = '  +6.64%  '
